$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.441378666666667
$ws.Range("H2").Value = 7.324135999999999
$ws.Range("I2").Value = 0.1119936059016048
$ws.Range("J2").Value = 0.1119936059016048
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.48183333333333
$ws.Range("N2").Value = 103.4455
$ws.Range("O2").Value = 0.2193934666657199
$ws.Range("P2").Value = 0.2193934666657199
$ws.Range("Q2").Value = 84.18321228755555
$ws.Range("R2").Value = 757.6489105879999
$ws.Range("S2").Value = 0.02457066544314751
$ws.Range("T2").Value = 0.0245706654431475

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.441378666666667
$ws.Range("H3").Value = 7.324135999999999
$ws.Range("I3").Value = 0.1119936059016048
$ws.Range("J3").Value = 0.1119936059016048
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 55.607043
$ws.Range("N3").Value = 166.821129
$ws.Range("O3").Value = 0.3538043298587107
$ws.Range("P3").Value = 0.3538043298587107
$ws.Range("Q3").Value = 135.757848496616
$ws.Range("R3").Value = 1221.820636469544
$ws.Range("S3").Value = 0.03962382268447784
$ws.Range("T3").Value = 0.03962382268447783

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.441378666666667
$ws.Range("H4").Value = 7.324135999999999
$ws.Range("I4").Value = 0.1119936059016048
$ws.Range("J4").Value = 0.1119936059016048
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 30.40497333333333
$ws.Range("N4").Value = 91.21491999999999
$ws.Range("O4").Value = 0.1934541136195998
$ws.Range("P4").Value = 0.1934541136195998
$ws.Range("Q4").Value = 74.23005325656888
$ws.Range("R4").Value = 668.0704793091198
$ws.Range("S4").Value = 0.02166562376075773
$ws.Range("T4").Value = 0.02166562376075773

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.441378666666667
$ws.Range("H5").Value = 7.324135999999999
$ws.Range("I5").Value = 0.1119936059016048
$ws.Range("J5").Value = 0.1119936059016048
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.91468833333333
$ws.Range("N5").Value = 89.744065
$ws.Range("O5").Value = 0.1903346354652808
$ws.Range("P5").Value = 0.1903346354652808
$ws.Range("Q5").Value = 73.03308191698223
$ws.Range("R5").Value = 657.29773725284
$ws.Range("S5").Value = 0.02131626215372427
$ws.Range("T5").Value = 0.02131626215372426

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.441378666666667
$ws.Range("H6").Value = 7.324135999999999
$ws.Range("I6").Value = 0.1119936059016048
$ws.Range("J6").Value = 0.1119936059016048
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.760378
$ws.Range("N6").Value = 20.281134
$ws.Range("O6").Value = 0.04301345439068881
$ws.Range("P6").Value = 0.04301345439068881
$ws.Range("Q6").Value = 16.50464262780267
$ws.Range("R6").Value = 148.541783650224
$ws.Range("S6").Value = 0.004817231859497455
$ws.Range("T6").Value = 0.004817231859497455

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.88577866666667
$ws.Range("H7").Value = 38.657336
$ws.Range("I7").Value = 0.5911106037886134
$ws.Range("J7").Value = 0.5911106037886134
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 34.48183333333333
$ws.Range("N7").Value = 103.4455
$ws.Range("O7").Value = 0.2193934666657199
$ws.Range("P7").Value = 0.2193934666657199
$ws.Range("Q7").Value = 444.3252723542222
$ws.Range("R7").Value = 3998.927451188
$ws.Range("S7").Value = 0.1296858045480507
$ws.Range("T7").Value = 0.1296858045480507

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.88577866666667
$ws.Range("H8").Value = 38.657336
$ws.Range("I8").Value = 0.5911106037886134
$ws.Range("J8").Value = 0.5911106037886134
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 55.607043
$ws.Range("N8").Value = 166.821129
$ws.Range("O8").Value = 0.3538043298587107
$ws.Range("P8").Value = 0.3538043298587107
$ws.Range("Q8").Value = 716.540048405816
$ws.Range("R8").Value = 6448.860435652344
$ws.Range("S8").Value = 0.2091374910458082
$ws.Range("T8").Value = 0.2091374910458082

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.88577866666667
$ws.Range("H9").Value = 38.657336
$ws.Range("I9").Value = 0.5911106037886134
$ws.Range("J9").Value = 0.5911106037886134
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.40497333333333
$ws.Range("N9").Value = 91.21491999999999
$ws.Range("O9").Value = 0.1934541136195998
$ws.Range("P9").Value = 0.1934541136195998
$ws.Range("Q9").Value = 391.7917567392356
$ws.Range("R9").Value = 3526.12581065312
$ws.Range("S9").Value = 0.1143527779070726
$ws.Range("T9").Value = 0.1143527779070726

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.88577866666667
$ws.Range("H10").Value = 38.657336
$ws.Range("I10").Value = 0.5911106037886134
$ws.Range("J10").Value = 0.5911106037886134
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 29.91468833333333
$ws.Range("N10").Value = 89.744065
$ws.Range("O10").Value = 0.1903346354652808
$ws.Range("P10").Value = 0.1903346354652808
$ws.Range("Q10").Value = 385.4740527456489
$ws.Range("R10").Value = 3469.26647471084
$ws.Range("S10").Value = 0.1125088212917677
$ws.Range("T10").Value = 0.1125088212917677

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 12.88577866666667
$ws.Range("H11").Value = 38.657336
$ws.Range("I11").Value = 0.5911106037886134
$ws.Range("J11").Value = 0.5911106037886134
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.760378
$ws.Range("N11").Value = 20.281134
$ws.Range("O11").Value = 0.04301345439068881
$ws.Range("P11").Value = 0.04301345439068881
$ws.Range("Q11").Value = 87.11273461100267
$ws.Range("R11").Value = 784.0146114990241
$ws.Range("S11").Value = 0.02542570899591404
$ws.Range("T11").Value = 0.02542570899591405

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.472111000000001
$ws.Range("H12").Value = 19.416333
$ws.Range("I12").Value = 0.2968957903097819
$ws.Range("J12").Value = 0.2968957903097818
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 34.48183333333333
$ws.Range("N12").Value = 103.4455
$ws.Range("O12").Value = 0.2193934666657199
$ws.Range("P12").Value = 0.2193934666657199
$ws.Range("Q12").Value = 223.1702528168334
$ws.Range("R12").Value = 2008.5322753515
$ws.Range("S12").Value = 0.0651369966745217
$ws.Range("T12").Value = 0.06513699667452169

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.472111000000001
$ws.Range("H13").Value = 19.416333
$ws.Range("I13").Value = 0.2968957903097819
$ws.Range("J13").Value = 0.2968957903097818
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 55.607043
$ws.Range("N13").Value = 166.821129
$ws.Range("O13").Value = 0.3538043298587107
$ws.Range("P13").Value = 0.3538043298587107
$ws.Range("Q13").Value = 359.8949546777731
$ws.Range("R13").Value = 3239.054592099957
$ws.Range("S13").Value = 0.1050430161284247
$ws.Range("T13").Value = 0.1050430161284247

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.472111000000001
$ws.Range("H14").Value = 19.416333
$ws.Range("I14").Value = 0.2968957903097819
$ws.Range("J14").Value = 0.2968957903097818
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.40497333333333
$ws.Range("N14").Value = 91.21491999999999
$ws.Range("O14").Value = 0.1934541136195998
$ws.Range("P14").Value = 0.1934541136195998
$ws.Range("Q14").Value = 196.7843623653733
$ws.Range("R14").Value = 1771.05926128836
$ws.Range("S14").Value = 0.05743571195176941
$ws.Range("T14").Value = 0.0574357119517694

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.472111000000001
$ws.Range("H15").Value = 19.416333
$ws.Range("I15").Value = 0.2968957903097819
$ws.Range("J15").Value = 0.2968957903097818
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 29.91468833333333
$ws.Range("N15").Value = 89.744065
$ws.Range("O15").Value = 0.1903346354652808
$ws.Range("P15").Value = 0.1903346354652808
$ws.Range("Q15").Value = 193.6111834237384
$ws.Range("R15").Value = 1742.500650813645
$ws.Range("S15").Value = 0.05650955201978877
$ws.Range("T15").Value = 0.05650955201978877

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.472111000000001
$ws.Range("H16").Value = 19.416333
$ws.Range("I16").Value = 0.2968957903097819
$ws.Range("J16").Value = 0.2968957903097818
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.760378
$ws.Range("N16").Value = 20.281134
$ws.Range("O16").Value = 0.04301345439068881
$ws.Range("P16").Value = 0.04301345439068881
$ws.Range("Q16").Value = 43.75391681795801
$ws.Range("R16").Value = 393.785251361622
$ws.Range("S16").Value = 0.01277051353527731
$ws.Range("T16").Value = 0.01277051353527731
